# Apply the diff: append one new data row to each of the 5 worksheets
# (status, neighbors, links, routes get a new row 3; topology gets new
# rows 4 and 5), reusing existing shared strings ("10.0.0.2", "10.0.0.1",
# "mesh0").

$wb = $excel.ActiveWorkbook

# --- status sheet: new row 3 ---
$ws = $wb.Worksheets.Item("status")
$ws.Cells.Item(3, 1).Value = 1465993371
$ws.Cells.Item(3, 2).Value = 1041489509

# --- neighbors sheet: new row 3 ---
$ws = $wb.Worksheets.Item("neighbors")
$ws.Cells.Item(3, 1).Value = 1465993371
$ws.Cells.Item(3, 2).Value = "10.0.0.2"
$ws.Cells.Item(3, 3).Value = $true
$ws.Cells.Item(3, 4).Value = $false
$ws.Cells.Item(3, 5).Value = $false
$ws.Cells.Item(3, 6).Value = 3
$ws.Cells.Item(3, 7).Value = 0

# --- links sheet: new row 3 ---
$ws = $wb.Worksheets.Item("links")
$ws.Cells.Item(3, 1).Value = 1465993371
$ws.Cells.Item(3, 2).Value = "10.0.0.1"
$ws.Cells.Item(3, 3).Value = "10.0.0.2"
$ws.Cells.Item(3, 4).Value = 38076
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1024

# --- routes sheet: new row 3 ---
$ws = $wb.Worksheets.Item("routes")
$ws.Cells.Item(3, 1).Value = 1465993371
$ws.Cells.Item(3, 2).Value = "10.0.0.2"
$ws.Cells.Item(3, 3).Value = 32
$ws.Cells.Item(3, 4).Value = "10.0.0.2"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 1024
$ws.Cells.Item(3, 7).Value = "mesh0"

# --- topology sheet: new rows 4 and 5 ---
$ws = $wb.Worksheets.Item("topology")
$ws.Cells.Item(4, 1).Value = 1465993371
$ws.Cells.Item(4, 2).Value = "10.0.0.2"
$ws.Cells.Item(4, 3).Value = "10.0.0.1"
$ws.Cells.Item(4, 4).Value = 1
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 1024
$ws.Cells.Item(4, 7).Value = 0

$ws.Cells.Item(5, 1).Value = 1465993371
$ws.Cells.Item(5, 2).Value = "10.0.0.1"
$ws.Cells.Item(5, 3).Value = "10.0.0.2"
$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 1024
$ws.Cells.Item(5, 7).Value = 283210
